$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.311.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.418.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.77%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  +6.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.585'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.20'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("E12").Value = '  +2.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '679.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.966.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.64'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.377.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.411.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.910'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.39%  '
$ws.Range("E22").Value = '  -2.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.90'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.41%  '
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '558.43'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.608.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.70%  '
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0734'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.77%  '
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.69'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.28%  '
$ws.Range("E43").Value = '  +5.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0424'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("E48").Value = '  +4.57%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("E51").Value = '  +2.46%  '
